$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pl_mw results for the 380 kV case (rows 2-25, columns B,D,E,F,G,I,J,L,M,N,O).
# Columns A, C, H, K are all zero in this case and are left untouched.
$newValues = @{
    2 = @{ "B" = 0.9633329015805998; "D" = 0.203555524993007; "E" = 0.207163166437752; "F" = 1.359360498981786; "G" = 0.002441103864676292; "I" = 0.9106198446675613; "J" = 0.2651511547096064; "L" = 0.3842297507223407; "M" = 0.3043433590389668; "N" = 1.648363412763231; "O" = 3.353242973395652 }
    3 = @{ "B" = 0.9125628844234939; "D" = 0.2033081500342391; "E" = 0.204701574945318; "F" = 1.358013640946311; "G" = 0.002444343894617368; "I" = 0.9232703879013169; "J" = 0.2594556388490616; "L" = 0.3580789927378021; "M" = 0.2863196614769095; "N" = 1.649826801869807; "O" = 3.327764594268018 }
    4 = @{ "B" = 0.8815858646239008; "D" = 0.2031702378077895; "E" = 0.2031911967146733; "F" = 1.357969042828422; "G" = 0.002446442261447968; "I" = 0.9314850100037857; "J" = 0.2559483652670309; "L" = 0.3420906157806058; "M" = 0.2753122552522953; "N" = 1.651291340439258; "O" = 3.31420766776597 }
    5 = @{ "B" = 0.8690128740341549; "D" = 0.2031175701993817; "E" = 0.2025760277623583; "F" = 1.358147588536418; "G" = 0.002447324848989358; "I" = 0.9349449483299654; "J" = 0.2545166808442758; "L" = 0.3355928317585608; "M" = 0.2708418931345022; "N" = 1.65203095495697; "O" = 3.309207635185828 }
    6 = @{ "B" = 0.8669282152500841; "D" = 0.2031090386874226; "E" = 0.2024739006288527; "F" = 1.358189116885114; "G" = 0.002447473064621677; "I" = 0.9355262572556526; "J" = 0.254278806894277; "L" = 0.3345149560961147; "M" = 0.2701005241390462; "N" = 1.652162406412231; "O" = 3.308409056593149 }
    7 = @{ "B" = 0.8814160953417343; "D" = 0.2031695131840401; "E" = 0.2031828989471123; "F" = 1.357970654270062; "G" = 0.002446454052901374; "I" = 0.9315312168166576; "J" = 0.2559290667928096; "L" = 0.3420029125146016; "M" = 0.2752519041877193; "N" = 1.651300736301295; "O" = 3.314138112125931 }
    8 = @{ "B" = 0.945787377241686; "D" = 0.2034673398583777; "E" = 0.2063142264021565; "F" = 1.35873370590091; "G" = 0.002442198459504382; "I" = 0.9148889359887606; "J" = 0.2631895311960903; "L" = 0.375199047988545; "M" = 0.2981167043529425; "N" = 1.64875079390869; "O" = 3.344024727452449 }
    9 = @{ "B" = 1.073531824984542; "D" = 0.2041615491344118; "E" = 0.2124607921522248; "F" = 1.366440666081928; "G" = 0.002434714120841377; "I" = 0.8858018912990158; "J" = 0.277341533814095; "L" = 0.4408226382452369; "M" = 0.3434111152000767; "N" = 1.648221954287834; "O" = 3.419210033938441 }
    10 = @{ "B" = 1.168260526646918; "D" = 0.2047378671560409; "E" = 0.2169777313161951; "F" = 1.375896032344471; "G" = 0.002429734869885792; "I" = 0.8665949904700181; "J" = 0.2876812782573097; "L" = 0.4893406864997019; "M" = 0.3769529009738477; "N" = 1.650535316171229; "O" = 3.484593135552245 }
    11 = @{ "B" = 1.211535637291888; "D" = 0.2050142633558991; "E" = 0.2190323039618036; "F" = 1.381022839609074; "G" = 0.002427581352575743; "I" = 0.8583270420409104; "J" = 0.2923714614502018; "L" = 0.5114757511659036; "M" = 0.3922664157455102; "N" = 1.652169543621355; "O" = 3.516549337823392 }
    12 = @{ "B" = 1.227948045228118; "D" = 0.2051209566455938; "E" = 0.2198102335443259; "F" = 1.383083016678697; "G" = 0.002426781829504029; "I" = 0.8552637059084454; "J" = 0.2941454683614921; "L" = 0.5198665469870036; "M" = 0.3980728677442826; "N" = 1.652871630485791; "O" = 3.528969049548095 }
    13 = @{ "B" = 1.224412243728068; "D" = 0.2050978884182797; "E" = 0.2196426976034971; "F" = 1.382634037587408; "G" = 0.002426953312123047; "I" = 0.8559204445677011; "J" = 0.2937634981469444; "L" = 0.5180590583078128; "M" = 0.3968220143834387; "N" = 1.652716728180678; "O" = 3.526280067408038 }
    14 = @{ "B" = 1.212885400468338; "D" = 0.2050230005332452; "E" = 0.2190963069101421; "F" = 1.381189951026343; "G" = 0.002427515255847527; "I" = 0.8580736652434346; "J" = 0.2925174521925555; "L" = 0.5121658947717549; "M" = 0.3927439667783617; "N" = 1.652225638917315; "O" = 3.517564727934428 }
    15 = @{ "B" = 1.205828104854106; "D" = 0.2049773931310455; "E" = 0.218761612924709; "F" = 1.380320874936729; "G" = 0.002427861539700079; "I" = 0.8594013747815055; "J" = 0.2917539403399587; "L" = 0.5085572845460149; "M" = 0.3902470176682371; "N" = 1.65193566022478; "O" = 3.512267831234738 }
    16 = @{ "B" = 1.165435968154554; "D" = 0.2047200887691929; "E" = 0.2168434506416936; "F" = 1.375577599892324; "G" = 0.002429877845759314; "I" = 0.8671447705061492; "J" = 0.2873744822705078; "L" = 0.487895354800969; "M" = 0.3759532048198508; "N" = 1.650440187198328; "O" = 3.482549285807494 }
    17 = @{ "B" = 1.140702609783261; "D" = 0.2045658719299794; "E" = 0.215666622719727; "F" = 1.372879229930177; "G" = 0.002431143305303445; "I" = 0.8720153531667698; "J" = 0.2846842942912957; "L" = 0.4752359771599686; "M" = 0.367198278201684; "N" = 1.649671484147021; "O" = 3.464884998802063 }
    18 = @{ "B" = 1.126493889389621; "D" = 0.2044785105732387; "E" = 0.2149897268221572; "F" = 1.371404896667642; "G" = 0.00243188166992017; "I" = 0.8748609587980578; "J" = 0.283135714314227; "L" = 0.4679606863418542; "M" = 0.3621678931835746; "N" = 1.649284139985866; "O" = 3.454933254834032 }
    19 = @{ "B" = 1.121686060294564; "D" = 0.2044491622208326; "E" = 0.2147605407414517; "F" = 1.370919055638936; "G" = 0.002432133474401495; "I" = 0.8758320172399738; "J" = 0.2826111802857127; "L" = 0.4654984510403324; "M" = 0.3604655988962548; "N" = 1.649162415998632; "O" = 3.451599528082596 }
    20 = @{ "B" = 1.143333740977027; "D" = 0.2045821500116034; "E" = 0.2157919001327002; "F" = 1.373158434163727; "G" = 0.002431007508219919; "I" = 0.8714922989002005; "J" = 0.2849708002551097; "L" = 0.4765829669470349; "M" = 0.3681297178765846; "N" = 1.649747646121853; "O" = 3.466743832866428 }
    21 = @{ "B" = 1.216270443002372; "D" = 0.2050449420198461; "E" = 0.2192567981121734; "F" = 1.381610890502714; "G" = 0.002427349766696458; "I" = 0.8574393781750178; "J" = 0.2928835031689943; "L" = 0.5138966266431737; "M" = 0.3939415861591513; "N" = 1.652367628024479; "O" = 3.520115985136442 }
    22 = @{ "B" = 1.264084329961747; "D" = 0.2053592162832061; "E" = 0.2215207410731672; "F" = 1.387827347898565; "G" = 0.002425052254181857; "I" = 0.8486487822407511; "J" = 0.2980428212446498; "L" = 0.5383339117312005; "M" = 0.4108550025281659; "N" = 1.654564930673217; "O" = 3.556854928287464 }
    23 = @{ "B" = 1.238552233678547; "D" = 0.2051904071015755; "E" = 0.2203125048377075; "F" = 1.384446143884517; "G" = 0.002426269992279939; "I" = 0.8533044353549979; "J" = 0.2952903466242134; "L" = 0.5252867932804293; "M" = 0.4018241051809142; "N" = 1.653347952906586; "O" = 3.537076617681976 }
    24 = @{ "B" = 1.142144172777023; "D" = 0.2045747866422545; "E" = 0.2157352632132188; "F" = 1.373031966094942; "G" = 0.002431068868149144; "I" = 0.8717286301866842; "J" = 0.2848412769900222; "L" = 0.4759739842498334; "M" = 0.3677086048153768; "N" = 1.649713043204287; "O" = 3.465902819521204 }
    25 = @{ "B" = 1.038816533331669; "D" = 0.2039620441496517; "E" = 0.2107976084361916; "F" = 1.363690126009004; "G" = 0.002436647222660006; "I" = 0.8932908129163302; "J" = 0.2735227871354837; "L" = 0.4230150890209075; "M" = 0.3604655988962548; "N" = 1.647888821900978; "O" = 3.39709170153742 }
}

foreach ($rowNum in $newValues.Keys) {
    $rowData = $newValues[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
